# Replace every "Non-bank financial services" sector label with the new
# "Financial Services" label (EGX re-categorized sectors), and update the
# active selection to C17, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldLabel = "Non-bank financial services"
$newLabel = "Financial Services"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

$changed = 0
for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -eq $oldLabel) {
        $cell.Value = $newLabel
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed cells from '$oldLabel' to '$newLabel'"

# Update the selection shown in the sheet view to C17, as in the target file.
$ws.Range("C17").Select()
